$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.117.33"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.654.72"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "'218.81"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'0.5257"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "'0.2682"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "'0.06377"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'20.58"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "'0.07685"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "'4.611"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "1.692.55"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "1.882.59"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "'0.5632"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "0.0₅8247"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "'65.68"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "26.103.83"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'4.696"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'10.34"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "'190.33"
$ws.Range("E22").Value = "  -4.69%  "
$ws.Range("D23").Value = "'5.991"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'146.61"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "'0.1202"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "'7.260"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "'15.99"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D29").Value = "'1.522"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "'0.05658"
$ws.Range("E30").Value = "  -4.25%  "
$ws.Range("D31").Value = "'1.277"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "'3.498"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "'3.383"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").Value = "'1.581"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "'2.794"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").Value = "'0.9490"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "'2.409"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").Value = "'0.5788"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'0.01596"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").Value = "'5.971"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'0.8371"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").Value = "1.023.56"
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("D44").Value = "'101.27"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "1.792.75"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "'58.32"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "'0.05332"
$ws.Range("E49").Value = "  +3.63%  "
$ws.Range("D50").Value = "'8.057"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'0.4340"
$ws.Range("E51").Value = "  -1.65%  "
